$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): rows 2-15, column F ("想去人数") updated
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 55
$ws1.Range("F3").Value = 11559
$ws1.Range("F4").Value = 207
$ws1.Range("F5").Value = 328
$ws1.Range("F6").Value = 222
$ws1.Range("F7").Value = 11514
$ws1.Range("F8").Value = 476
$ws1.Range("F9").Value = 1163
$ws1.Range("F10").Value = 83
$ws1.Range("F11").Value = 1756
$ws1.Range("F12").Value = 5714
$ws1.Range("F14").Value = 3500
$ws1.Range("F15").Value = 179

# Sheet "全部类型" (sheet4): rows 3-18, column F ("想去人数") updated
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 55
$ws4.Range("F5").Value = 11559
$ws4.Range("F6").Value = 207
$ws4.Range("F7").Value = 328
$ws4.Range("F8").Value = 222
$ws4.Range("F9").Value = 11515
$ws4.Range("F10").Value = 476
$ws4.Range("F11").Value = 1163
$ws4.Range("F12").Value = 83
$ws4.Range("F13").Value = 1756
$ws4.Range("F15").Value = 5715
$ws4.Range("F17").Value = 3500
$ws4.Range("F18").Value = 179
